# Apply scheduled-runner updates to Sheets (leve profit data refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("H112").Value = 34091908
$ws.Range("J112").Value = 34091908
$ws.Range("L112").Value = 102275724
$ws.Range("N112").Value = -102277940
$ws.Range("H138").Value = 6946996
$ws.Range("I138").Value = 2676.238
$ws.Range("J138").Value = 9806421
$ws.Range("K138").Value = 8028.714
$ws.Range("L138").Value = 29419263
$ws.Range("M138").Value = -2888.714
$ws.Range("N138").Value = -29429543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2721.9487
$ws.Range("I32").Value = 2218.1785
$ws.Range("J32").Value = 4004.2727
$ws.Range("K32").Value = 2218.1785
$ws.Range("L32").Value = 4004.2727
$ws.Range("M32").Value = -1931.1785
$ws.Range("N32").Value = -4578.2727
$ws.Range("H60").Value = 49775
$ws.Range("I60").Value = 49775
$ws.Range("K60").Value = 49775
$ws.Range("M60").Value = -49042

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11948.467
$ws.Range("I86").Value = 9308
$ws.Range("J86").Value = 14258.875
$ws.Range("K86").Value = 9308
$ws.Range("L86").Value = 14258.875
$ws.Range("M86").Value = -8185
$ws.Range("N86").Value = -16504.875
$ws.Range("H89").Value = 11948.467
$ws.Range("I89").Value = 9308
$ws.Range("J89").Value = 14258.875
$ws.Range("K89").Value = 46540
$ws.Range("L89").Value = 71294.375
$ws.Range("M89").Value = -40924
$ws.Range("N89").Value = -82526.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1581.2982
$ws.Range("I31").Value = 1020.91174
$ws.Range("J31").Value = 2409.6956
$ws.Range("K31").Value = 1020.91174
$ws.Range("L31").Value = 2409.6956
$ws.Range("M31").Value = -725.91174
$ws.Range("N31").Value = -2999.6956
$ws.Range("H34").Value = 1581.2982
$ws.Range("I34").Value = 1020.91174
$ws.Range("J34").Value = 2409.6956
$ws.Range("K34").Value = 1020.91174
$ws.Range("L34").Value = 2409.6956
$ws.Range("M34").Value = -818.91174
$ws.Range("N34").Value = -2813.6956
$ws.Range("H44").Value = 6000
$ws.Range("J44").Value = 6000
$ws.Range("L44").Value = 6000
$ws.Range("N44").Value = -6884
$ws.Range("H45").Value = 6499.8335
$ws.Range("I45").Value = 1999
$ws.Range("K45").Value = 1999
$ws.Range("M45").Value = -1406
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 845.913
$ws.Range("I122").Value = 569.4545000000001
$ws.Range("K122").Value = 5125.0905
$ws.Range("M122").Value = -2675.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2431.818
$ws.Range("I80").Value = 2414.2856
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 2414.2856
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -1416.2856
$ws.Range("N80").Value = -4796
$ws.Range("H83").Value = 2431.818
$ws.Range("I83").Value = 2414.2856
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 12071.428
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -7079.428
$ws.Range("N83").Value = -23984
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 2417.7
$ws.Range("I126").Value = 2078
$ws.Range("J126").Value = 2502.625
$ws.Range("K126").Value = 6234
$ws.Range("L126").Value = 7507.875
$ws.Range("M126").Value = -3764
$ws.Range("N126").Value = -12447.875
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 70996
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 70996
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 70996
$ws.Range("N128").Value = -80956
$ws.Range("H129").Value = 49999.8
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49999.8
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49999.8
$ws.Range("N129").Value = -59999.8
$ws.Range("H130").Value = 44875
$ws.Range("I130").Value = 30500
$ws.Range("J130").Value = 49666.668
$ws.Range("K130").Value = 30500
$ws.Range("L130").Value = 49666.668
$ws.Range("M130").Value = -25480
$ws.Range("N130").Value = -59706.668
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 2484.2144
$ws.Range("I132").Value = 2123.138
$ws.Range("J132").Value = 3289.6924
$ws.Range("K132").Value = 6369.414
$ws.Range("L132").Value = 9869.0772
$ws.Range("M132").Value = -3839.414
$ws.Range("N132").Value = -14929.0772
$ws.Range("H133").Value = 16563.158
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 16563.158
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 16563.158
$ws.Range("N133").Value = -26683.158
$ws.Range("H134").Value = 24206
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 24206
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 72618
$ws.Range("N134").Value = -77688
$ws.Range("H135").Value = 1000000000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1000000000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 1000000000
$ws.Range("N135").Value = -1000010140
$ws.Range("H136").Value = 32663
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 32663
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 97989
$ws.Range("N136").Value = -103089
$ws.Range("H137").Value = 60000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 60000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200
$ws.Range("H138").Value = 62666.668
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 62666.668
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 62666.668
$ws.Range("N138").Value = -72946.66800000001
$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
$ws.Range("H140").Value = 56625
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 56625
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 56625
$ws.Range("N140").Value = -66985
$ws.Range("H141").Value = 85000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 85000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 85000
$ws.Range("N141").Value = -95360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 34494.5
$ws.Range("J48").Value = 18989
$ws.Range("L48").Value = 18989
$ws.Range("N48").Value = -20311
$ws.Range("H82").Value = 43207.668
$ws.Range("I82").Value = 167999.67
$ws.Range("J82").Value = 1610.3334
$ws.Range("K82").Value = 167999.67
$ws.Range("L82").Value = 1610.3334
$ws.Range("M82").Value = -167638.67
$ws.Range("N82").Value = -2332.3334
$ws.Range("H85").Value = 43207.668
$ws.Range("I85").Value = 167999.67
$ws.Range("J85").Value = 1610.3334
$ws.Range("K85").Value = 167999.67
$ws.Range("L85").Value = 1610.3334
$ws.Range("M85").Value = -166751.67
$ws.Range("N85").Value = -4106.3334
$ws.Range("H122").Value = 3382.1072
$ws.Range("I122").Value = 2522.111
$ws.Range("K122").Value = 7566.333
$ws.Range("M122").Value = -5116.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7378.9473
$ws.Range("I54").Value = 7200
$ws.Range("J54").Value = 7400
$ws.Range("K54").Value = 7200
$ws.Range("L54").Value = 7400
$ws.Range("M54").Value = -6680
$ws.Range("N54").Value = -8440
$ws.Range("H132").Value = 16131623
$ws.Range("I132").Value = 25002038
$ws.Range("K132").Value = 75006114
$ws.Range("M132").Value = -75003584
